$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 2292.1072
$ws.Range("I15").Value = 2292.1072
$ws.Range("K15").Value = 6876.321599999999
$ws.Range("M15").Value = -6707.321599999999

$ws.Range("H19").Value = 2078.1667
$ws.Range("I19").Value = 1868.25
$ws.Range("J19").Value = 2498
$ws.Range("K19").Value = 1868.25
$ws.Range("L19").Value = 2498
$ws.Range("M19").Value = -1693.25
$ws.Range("N19").Value = -2848

$ws.Range("H107").Value = 2297.7144
$ws.Range("I107").Value = 1096.8
$ws.Range("J107").Value = 5300
$ws.Range("K107").Value = 1096.8
$ws.Range("L107").Value = 5300
$ws.Range("M107").Value = 823.2
$ws.Range("N107").Value = -9140

$ws.Range("H125").Value = 3565.9285
$ws.Range("I125").Value = 1993.1818
$ws.Range("J125").Value = 9332.666999999999
$ws.Range("K125").Value = 17938.6362
$ws.Range("L125").Value = 83994.003
$ws.Range("M125").Value = -15478.6362
$ws.Range("N125").Value = -88914.003

$ws.Range("H138").Value = 3081.258
$ws.Range("I138").Value = 1654.4166
$ws.Range("J138").Value = 3982.4211
$ws.Range("K138").Value = 4963.2498
$ws.Range("L138").Value = 11947.2633
$ws.Range("M138").Value = 176.7502000000004
$ws.Range("N138").Value = -22227.2633

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 3029.0625
$ws.Range("I74").Value = 3029.0625
$ws.Range("K74").Value = 3029.0625
$ws.Range("M74").Value = -2155.0625

$ws.Range("H77").Value = 3029.0625
$ws.Range("I77").Value = 3029.0625
$ws.Range("K77").Value = 15145.3125
$ws.Range("M77").Value = -10777.3125

$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2765.037
$ws.Range("I20").Value = 2119
$ws.Range("J20").Value = 3863.3
$ws.Range("K20").Value = 2119
$ws.Range("L20").Value = 3863.3
$ws.Range("M20").Value = -1872
$ws.Range("N20").Value = -4357.3

$ws.Range("H94").Value = 1661.75
$ws.Range("I94").Value = 1753.7727
$ws.Range("J94").Value = 649.5
$ws.Range("K94").Value = 1753.7727
$ws.Range("L94").Value = 649.5
$ws.Range("M94").Value = -1302.7727
$ws.Range("N94").Value = -1551.5

$ws.Range("H134").Value = 4650.25
$ws.Range("I134").Value = 4650.25
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 13950.75
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -11415.75
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2125.8936
$ws.Range("I31").Value = 1624.8611
$ws.Range("J31").Value = 3765.6365
$ws.Range("K31").Value = 1624.8611
$ws.Range("L31").Value = 3765.6365
$ws.Range("M31").Value = -1329.8611
$ws.Range("N31").Value = -4355.636500000001

$ws.Range("H34").Value = 2125.8936
$ws.Range("I34").Value = 1624.8611
$ws.Range("J34").Value = 3765.6365
$ws.Range("K34").Value = 1624.8611
$ws.Range("L34").Value = 3765.6365
$ws.Range("M34").Value = -1422.8611
$ws.Range("N34").Value = -4169.636500000001

$ws.Range("H51").Value = 24800
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 24800
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 24800
$ws.Range("N51").Value = -26272
$ws.Range("M51").ClearContents()

$ws.Range("H61").Value = 24800
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 24800
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 24800
$ws.Range("N61").Value = -25496
$ws.Range("M61").ClearContents()

$ws.Range("H86").Value = 11930.25
$ws.Range("I86").Value = 12277.667
$ws.Range("K86").Value = 12277.667
$ws.Range("M86").Value = -11154.667

$ws.Range("H89").Value = 11930.25
$ws.Range("I89").Value = 12277.667
$ws.Range("K89").Value = 61388.335
$ws.Range("M89").Value = -55772.335

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1844.9166
$ws.Range("I4").Value = 1923.2
$ws.Range("J4").Value = 1714.4445
$ws.Range("K4").Value = 5769.6
$ws.Range("L4").Value = 5143.333500000001
$ws.Range("M4").Value = -5657.6
$ws.Range("N4").Value = -5367.333500000001

$ws.Range("H6").Value = 10270.667
$ws.Range("I6").Value = 108.333336
$ws.Range("K6").Value = 325.000008
$ws.Range("M6").Value = -212.000008

$ws.Range("H7").Value = 353.27274
$ws.Range("I7").Value = 298.33334
$ws.Range("K7").Value = 895.0000200000001
$ws.Range("M7").Value = -783.0000200000001

$ws.Range("H17").Value = 1066
$ws.Range("I17").Value = 399
$ws.Range("K17").Value = 1197
$ws.Range("M17").Value = -1028

$ws.Range("H26").Value = 176791.83
$ws.Range("I26").Value = 252437.75
$ws.Range("J26").Value = 25500
$ws.Range("K26").Value = 757313.25
$ws.Range("L26").Value = 76500
$ws.Range("M26").Value = -757025.25
$ws.Range("N26").Value = -77076

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 1593.125
$ws.Range("I5").Value = 1463.5714
$ws.Range("K5").Value = 1463.5714
$ws.Range("M5").Value = -1351.5714

$ws.Range("H122").Value = 4239.5
$ws.Range("I122").Value = 1950
$ws.Range("J122").Value = 5002.6665
$ws.Range("K122").Value = 5850
$ws.Range("L122").Value = 15007.9995
$ws.Range("M122").Value = -3400
$ws.Range("N122").Value = -19907.9995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 3499
$ws.Range("I16").Value = 3499
$ws.Range("K16").Value = 3499
$ws.Range("M16").Value = -3329

$ws.Range("H22").Value = 978.4167
$ws.Range("J22").Value = 750
$ws.Range("L22").Value = 750
$ws.Range("N22").Value = -1340

$ws.Range("H27").Value = 978.4167
$ws.Range("J27").Value = 750
$ws.Range("L27").Value = 750
$ws.Range("N27").Value = -964

$ws.Range("H40").Value = 2000
$ws.Range("I40").Value = 2000
$ws.Range("K40").Value = 2000
$ws.Range("M40").Value = -1864

$ws.Range("H46").Value = 3759.8
$ws.Range("I46").Value = 2937.5
$ws.Range("K46").Value = 2937.5
$ws.Range("M46").Value = -2749.5

$ws.Range("H55").Value = 248.57895
$ws.Range("I55").Value = 199.66667
$ws.Range("K55").Value = 199.66667
$ws.Range("M55").Value = -26.66667000000001

$ws.Range("H122").Value = 3704.2
$ws.Range("I122").Value = 3504
$ws.Range("K122").Value = 10512
$ws.Range("M122").Value = -8062

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 447
$ws.Range("I107").Value = 263
$ws.Range("K107").Value = 789
$ws.Range("M107").Value = 1131

$ws.Range("H136").Value = 9724.691999999999
$ws.Range("I136").Value = 9562.200000000001
$ws.Range("J136").Value = 10266.333
$ws.Range("K136").Value = 28686.6
$ws.Range("L136").Value = 30798.999
$ws.Range("M136").Value = -26136.6
$ws.Range("N136").Value = -35898.999
